$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily clear the "Per Week Avg" formula cell's content; leaving it
# in place while the rows below it get pushed down by the upcoming insert
# causes the engine to register a spurious, unused number-format style.
$ws.Range("B30").ClearContents()

# Insert a blank row at 28. The new row 28 inherits the plain
# (non-highlighted) formatting of row 27 above it; the row that used to
# be row 28 (still carrying its highlighted look) simply shifts down to
# become row 29, along with the "Total"/"Per Week Avg" rows moving to
# 30/31.
$ws.Rows("28").Insert()

# New row 28 keeps the same date it always had, just with the lower
# (no-longer-highlighted) hours figure.
$ws.Range("A28").Value = [DateTime]"2020-01-27"
$ws.Range("B28").Value = 11.75

# Row 29 (the former row 28) is now the latest, highlighted entry.
$ws.Range("A29").Value = [DateTime]"2020-02-03"
$ws.Range("B29").Value = 11.75

# Fix up the formulas in the Total / Per Week Avg rows that got pushed down
$ws.Range("B30").Formula = "=SUM(B2:B29)"
$ws.Range("B31").Formula = "=B30/COUNT(B2:B29)"

# Update the selection to match the new active cell
$ws.Range("A28:B28").Select()
